$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 22
$ws.Range("H22").Value = 1980
$ws.Range("I22").Value = 1980
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 5940
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -5768
$ws.Range("N22").ClearContents()

# Row 33
$ws.Range("H33").Value = 216.33333
$ws.Range("I33").Value = 145.81818
$ws.Range("K33").Value = 145.81818
$ws.Range("M33").Value = 83.18181999999999

# Row 40
$ws.Range("H40").Value = 2566.7917
$ws.Range("I40").Value = 2464.6667
$ws.Range("J40").Value = 2737
$ws.Range("K40").Value = 2464.6667
$ws.Range("L40").Value = 2737
$ws.Range("M40").Value = -2289.6667
$ws.Range("N40").Value = -3087

# Row 64
$ws.Range("H64").Value = 3502.0833
$ws.Range("I64").Value = 3248.718
$ws.Range("J64").Value = 4600
$ws.Range("K64").Value = 3248.718
$ws.Range("L64").Value = 4600
$ws.Range("M64").Value = -3000.718
$ws.Range("N64").Value = -5096

# Row 67
$ws.Range("H67").Value = 3502.0833
$ws.Range("I67").Value = 3248.718
$ws.Range("J67").Value = 4600
$ws.Range("K67").Value = 3248.718
$ws.Range("L67").Value = 4600
$ws.Range("M67").Value = -2390.718
$ws.Range("N67").Value = -6316

# Row 76
$ws.Range("H76").Value = 3391.3044
$ws.Range("I76").Value = 2999.4443
$ws.Range("J76").Value = 4802
$ws.Range("K76").Value = 2999.4443
$ws.Range("L76").Value = 4802
$ws.Range("M76").Value = -2684.4443
$ws.Range("N76").Value = -5432

# Row 79
$ws.Range("H79").Value = 3391.3044
$ws.Range("I79").Value = 2999.4443
$ws.Range("J79").Value = 4802
$ws.Range("K79").Value = 2999.4443
$ws.Range("L79").Value = 4802
$ws.Range("M79").Value = -1907.4443
$ws.Range("N79").Value = -6986

$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Range("H102").Value = 2791
$ws.Range("I102").Value = 1933.3334
$ws.Range("J102").Value = 4077.5
$ws.Range("K102").Value = 1933.3334
$ws.Range("L102").Value = 4077.5
$ws.Range("M102").Value = -311.3334
$ws.Range("N102").Value = -7321.5

$ws = $wb.Worksheets.Item("BSM")
# Row 34
$ws.Range("H34").Value = 19830
$ws.Range("J34").Value = 19830
$ws.Range("L34").Value = 19830
$ws.Range("N34").Value = -20058

# Row 99
$ws.Range("H99").Value = 1806.1765
$ws.Range("I99").Value = 1342.5
$ws.Range("J99").Value = 3970
$ws.Range("K99").Value = 1342.5
$ws.Range("L99").Value = 3970
$ws.Range("M99").Value = 155.5
$ws.Range("N99").Value = -6966

# Row 105
$ws.Range("H105").Value = 2382.2942
$ws.Range("I105").Value = 2140.75
$ws.Range("J105").Value = 2962
$ws.Range("K105").Value = 2140.75
$ws.Range("L105").Value = 2962
$ws.Range("M105").Value = -393.75
$ws.Range("N105").Value = -6456

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 4715
$ws.Range("I62").Value = 4082
$ws.Range("J62").Value = 5066.6665
$ws.Range("K62").Value = 4082
$ws.Range("L62").Value = 5066.6665
$ws.Range("M62").Value = -3458
$ws.Range("N62").Value = -6314.6665

# Row 65
$ws.Range("H65").Value = 4715
$ws.Range("I65").Value = 4082
$ws.Range("J65").Value = 5066.6665
$ws.Range("K65").Value = 20410
$ws.Range("L65").Value = 25333.3325
$ws.Range("M65").Value = -17290
$ws.Range("N65").Value = -31573.3325

# Row 87
$ws.Range("H87").Value = 25330
$ws.Range("J87").Value = 25330
$ws.Range("L87").Value = 25330
$ws.Range("N87").Value = -27702

# Row 90
$ws.Range("H90").Value = 25330
$ws.Range("J90").Value = 25330
$ws.Range("L90").Value = 75990
$ws.Range("N90").Value = -87846

$ws = $wb.Worksheets.Item("CUL")
# Row 48
$ws.Range("H48").Value = 2741.111
$ws.Range("I48").Value = 1500
$ws.Range("K48").Value = 4500
$ws.Range("M48").Value = -4250

# Row 82
$ws.Range("H82").Value = 104490.7
$ws.Range("I82").Value = 971
$ws.Range("J82").Value = 148856.28
$ws.Range("K82").Value = 2913
$ws.Range("L82").Value = 446568.84
$ws.Range("M82").Value = -2507
$ws.Range("N82").Value = -447380.84

# Row 85
$ws.Range("H85").Value = 104490.7
$ws.Range("I85").Value = 971
$ws.Range("J85").Value = 148856.28
$ws.Range("K85").Value = 2913
$ws.Range("L85").Value = 446568.84
$ws.Range("M85").Value = -1509
$ws.Range("N85").Value = -449376.84

# Row 88
$ws.Range("H88").Value = 1825
$ws.Range("J88").Value = 1825
$ws.Range("L88").Value = 5475
$ws.Range("N88").Value = -6331

# Row 91
$ws.Range("H91").Value = 1825
$ws.Range("J91").Value = 1825
$ws.Range("L91").Value = 5475
$ws.Range("N91").Value = -8439

# Row 92
$ws.Range("H92").Value = 875
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 1250
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 3750
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -6246

# Row 131
$ws.Range("H131").Value = 778.79
$ws.Range("I131").Value = 353.07693
$ws.Range("J131").Value = 842.4023
$ws.Range("K131").Value = 1059.23079
$ws.Range("L131").Value = 2527.2069
$ws.Range("M131").Value = 3980.76921
$ws.Range("N131").Value = -12607.2069

# Row 140
$ws.Range("H140").Value = 2429.2058
$ws.Range("I140").Value = 1565.7142
$ws.Range("J140").Value = 3824.077
$ws.Range("K140").Value = 4697.142599999999
$ws.Range("L140").Value = 11472.231
$ws.Range("M140").Value = 482.8574000000008
$ws.Range("N140").Value = -21832.231

$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 137500.08
$ws.Range("I21").Value = 4525.75
$ws.Range("J21").Value = 403448.75
$ws.Range("K21").Value = 4525.75
$ws.Range("L21").Value = 403448.75
$ws.Range("M21").Value = -4352.75
$ws.Range("N21").Value = -403794.75

# Row 24
$ws.Range("H24").Value = 50000
$ws.Range("J24").Value = 50000
$ws.Range("L24").Value = 50000
$ws.Range("N24").Value = -50346

# Row 30
$ws.Range("H30").Value = 137500.08
$ws.Range("I30").Value = 4525.75
$ws.Range("J30").Value = 403448.75
$ws.Range("K30").Value = 4525.75
$ws.Range("L30").Value = 403448.75
$ws.Range("M30").Value = -4420.75
$ws.Range("N30").Value = -403658.75

# Row 58
$ws.Range("H58").Value = 12000
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

# Row 70
$ws.Range("H70").Value = 5711.4414
$ws.Range("I70").Value = 4775.1875
$ws.Range("J70").Value = 6543.6665
$ws.Range("K70").Value = 4775.1875
$ws.Range("L70").Value = 6543.6665
$ws.Range("M70").Value = -4505.1875
$ws.Range("N70").Value = -7083.6665

# Row 73
$ws.Range("H73").Value = 5711.4414
$ws.Range("I73").Value = 4775.1875
$ws.Range("J73").Value = 6543.6665
$ws.Range("K73").Value = 4775.1875
$ws.Range("L73").Value = 6543.6665
$ws.Range("M73").Value = -3839.1875
$ws.Range("N73").Value = -8415.666499999999

$ws = $wb.Worksheets.Item("LTW")
# Row 44
$ws.Range("H44").Value = 16666.666
$ws.Range("J44").Value = 16666.666
$ws.Range("L44").Value = 16666.666
$ws.Range("N44").Value = -17578.666

# Row 46
$ws.Range("H46").Value = 1092.1428
$ws.Range("I46").Value = 911.25
$ws.Range("J46").Value = 1333.3334
$ws.Range("K46").Value = 911.25
$ws.Range("L46").Value = 1333.3334
$ws.Range("M46").Value = -723.25
$ws.Range("N46").Value = -1709.3334

# Row 100
$ws.Range("H100").Value = 2246
$ws.Range("I100").Value = 1277.4
$ws.Range("J100").Value = 5474.6665
$ws.Range("K100").Value = 1277.4
$ws.Range("L100").Value = 5474.6665
$ws.Range("M100").Value = -736.4000000000001
$ws.Range("N100").Value = -6556.6665

# Row 102
$ws.Range("H102").Value = 44000
$ws.Range("J102").Value = 44000
$ws.Range("L102").Value = 44000
$ws.Range("N102").Value = -50490

$ws = $wb.Worksheets.Item("WVR")
# Row 86
$ws.Range("H86").Value = 20226.818
$ws.Range("J86").Value = 20226.818
$ws.Range("L86").Value = 20226.818
$ws.Range("N86").Value = -22472.818

# Row 89
$ws.Range("H89").Value = 20226.818
$ws.Range("J89").Value = 20226.818
$ws.Range("L89").Value = 101134.09
$ws.Range("N89").Value = -112366.09

# Row 109
$ws.Range("H109").Value = 36263.332
$ws.Range("J109").Value = 36263.332
$ws.Range("L109").Value = 36263.332
$ws.Range("N109").Value = -39037.332
